# Generate Report for Archive
#
# The localization status for "957025f1-4f29-44bf-9392-00bd8277d1f7.md"
# moved from "Ready for handoff" back to "In Translation" (e.g. a new
# handoff cycle was kicked off for this file), so the generated report
# needs to reflect that on the Overview sheet as well as each per-locale
# status sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: columns are File Name | zh-cn | de-de | Latest Handoff Date
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B5").Value = "In Translation"
$overview.Range("C5").Value = "In Translation"

# Per-locale status sheets: columns include Source File Name | File Extension | Status | ...
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C5").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C5").Value = "In Translation"
